$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new column L for Russia data (shifts old L->M, old M->N)
$ws.Columns("L").Insert()

# Header row
$ws.Range("L1").Value = "Russia"

# Row labels (column A) - updated order for rows 6-9
$ws.Range("A2").Value = "Minimum tax of 2% on billionaires'`nwealth, in voluntary countries"
$ws.Range("A3").Value = "Bridgetown initiative: MDBs expanding sustainable`ninvestments in LICs, and at lower interest rates"
$ws.Range("A4").Value = "L&D: Developed countries financing a fund to help`nvulnerable countries cope with climate Loss and damage"
$ws.Range("A5").Value = "Expand Security Council to new permanent members (e.g.`nIndia, Brazil, African Union), restrict veto use"
$ws.Range("A6").Value = "Debt relief for vulnerable countries, suspending`npayments until they are more able to repay"
$ws.Range("A7").Value = "International levy on shipping carbon emissions,`nreturned to countries based on population"
$ws.Range("A8").Value = "At least 0.7% of developed countries' GDP in foreign aid"
$ws.Range("A9").Value = "Raise global minimum tax on profit from 15% to 35%,`nallocating revenues to countries based on sales"
$ws.Range("A10").Value = "NCQG: Developing countries providing `$300 bn a`nyear in climate finance for developing countries"
$ws.Range("A11").Value = "International levy on aviation carbon emissions, raising`nprices by 30%, returned to countries based on population"

# Data values
$ws.Range("B2").Value = 0.797728501892915
$ws.Range("C2").Value = 0.826012058570198
$ws.Range("D2").Value = 0.85842127450794
$ws.Range("E2").Value = 0.822487754215921
$ws.Range("F2").Value = 0.863876775817389
$ws.Range("G2").Value = 0.771221947124432
$ws.Range("H2").Value = 0.820473479187211
$ws.Range("I2").Value = 0.833078388530814
$ws.Range("J2").Value = 0.765444555149747
$ws.Range("K2").Value = 0.772249782420811
$ws.Range("L2").Value = 0.793823059348919
$ws.Range("M2").Value = 0.808205588588327
$ws.Range("N2").Value = 0.758763085765361
$ws.Range("B3").Value = 0.774794198353587
$ws.Range("C3").Value = 0.795595978937291
$ws.Range("D3").Value = 0.793261299509476
$ws.Range("E3").Value = 0.784892930783755
$ws.Range("F3").Value = 0.864422407284945
$ws.Range("G3").Value = 0.697208929057947
$ws.Range("H3").Value = 0.819169469417425
$ws.Range("I3").Value = 0.821194161810209
$ws.Range("J3").Value = 0.745648172295692
$ws.Range("K3").Value = 0.732817501320784
$ws.Range("L3").Value = 0.817116950550607
$ws.Range("M3").Value = 0.852437001452164
$ws.Range("N3").Value = 0.723005171748347
$ws.Range("B4").Value = 0.735490009514748
$ws.Range("C4").Value = 0.736258660508083
$ws.Range("D4").Value = 0.713556455527246
$ws.Range("E4").Value = 0.720391440962945
$ws.Range("F4").Value = 0.827442753819299
$ws.Range("G4").Value = 0.708068790596872
$ws.Range("H4").Value = 0.777216229032656
$ws.Range("I4").Value = 0.724437174674906
$ws.Range("J4").Value = 0.646238489843471
$ws.Range("K4").Value = 0.67903251711971
$ws.Range("L4").Value = 0.851034424656702
$ws.Range("M4").Value = 0.851960028038753
$ws.Range("N4").Value = 0.678845245140098
$ws.Range("B5").Value = 0.700619834710744
$ws.Range("C5").Value = 0.748409202153696
$ws.Range("D5").Value = 0.717020759449205
$ws.Range("E5").Value = 0.750158932763184
$ws.Range("F5").Value = 0.800661655988323
$ws.Range("G5").Value = 0.717479689795248
$ws.Range("H5").Value = 0.775217309431383
$ws.Range("I5").Value = 0.764109925653585
$ws.Range("J5").Value = 0.687684263344843
$ws.Range("K5").Value = 0.623577294273261
$ws.Range("L5").Value = 0.543297585547114
$ws.Range("M5").Value = 0.78738583872611
$ws.Range("N5").Value = 0.660738583552298
$ws.Range("B6").Value = 0.696536881236549
$ws.Range("C6").Value = 0.691183404054691
$ws.Range("D6").Value = 0.630233279970172
$ws.Range("E6").Value = 0.603125941881532
$ws.Range("F6").Value = 0.801138519399549
$ws.Range("G6").Value = 0.768633155276535
$ws.Range("H6").Value = 0.732948425350415
$ws.Range("I6").Value = 0.712952434354826
$ws.Range("J6").Value = 0.6282335145465
$ws.Range("K6").Value = 0.649007244283415
$ws.Range("L6").Value = 0.73705279435081
$ws.Range("M6").Value = 0.851313151457241
$ws.Range("N6").Value = 0.64805959650704
$ws.Range("B7").Value = 0.69309150193601
$ws.Range("C7").Value = 0.722995679308689
$ws.Range("D7").Value = 0.76683697471266
$ws.Range("E7").Value = 0.694591801308621
$ws.Range("F7").Value = 0.767042786406479
$ws.Range("G7").Value = 0.620477168841821
$ws.Range("H7").Value = 0.754737535106032
$ws.Range("I7").Value = 0.739810708802606
$ws.Range("J7").Value = 0.696830160275269
$ws.Range("K7").Value = 0.550584063687184
$ws.Range("L7").Value = 0.723263466239051
$ws.Range("M7").Value = 0.791664053032577
$ws.Range("N7").Value = 0.663301337003796
$ws.Range("B8").Value = 0.691438763376932
$ws.Range("C8").Value = 0.67858845970434
$ws.Range("D8").Value = 0.649566333705771
$ws.Range("E8").Value = 0.666863760631278
$ws.Range("F8").Value = 0.767588288389665
$ws.Range("G8").Value = 0.56800325399601
$ws.Range("H8").Value = 0.782834042853143
$ws.Range("I8").Value = 0.662139020074036
$ws.Range("J8").Value = 0.62415407154462
$ws.Range("K8").Value = 0.598307433828746
$ws.Range("L8").Value = 0.809636950891532
$ws.Range("M8").Value = 0.839647773958159
$ws.Range("N8").Value = 0.654362878650455
$ws.Range("B9").Value = 0.681889149102264
$ws.Range("C9").Value = 0.724314272431427
$ws.Range("D9").Value = 0.74804344232966
$ws.Range("E9").Value = 0.721006097532411
$ws.Range("F9").Value = 0.82158533790241
$ws.Range("G9").Value = 0.655850174362352
$ws.Range("H9").Value = 0.720464698516305
$ws.Range("I9").Value = 0.727314051621683
$ws.Range("J9").Value = 0.588476192495
$ws.Range("K9").Value = 0.66417147238743
$ws.Range("L9").Value = 0.479239298690242
$ws.Range("M9").Value = 0.764169437385325
$ws.Range("N9").Value = 0.648334728500467
$ws.Range("B10").Value = 0.677425665101721
$ws.Range("C10").Value = 0.684210526315789
$ws.Range("D10").Value = 0.666317918158563
$ws.Range("E10").Value = 0.68661094961858
$ws.Range("F10").Value = 0.751421893752358
$ws.Range("G10").Value = 0.63291031803108
$ws.Range("H10").Value = 0.744479707700628
$ws.Range("I10").Value = 0.663027150487923
$ws.Range("J10").Value = 0.629756686812235
$ws.Range("K10").Value = 0.562548331362135
$ws.Range("L10").Value = 0.8633755152337
$ws.Range("M10").Value = 0.821278106814353
$ws.Range("N10").Value = 0.607054320066253
$ws.Range("B11").Value = 0.547600468201327
$ws.Range("C11").Value = 0.553318077803204
$ws.Range("D11").Value = 0.602403959821429
$ws.Range("E11").Value = 0.540843569824836
$ws.Range("F11").Value = 0.561226840020469
$ws.Range("G11").Value = 0.53495446287913
$ws.Range("H11").Value = 0.567205056430566
$ws.Range("I11").Value = 0.553188855279033
$ws.Range("J11").Value = 0.505993094467125
$ws.Range("K11").Value = 0.45836014933668
$ws.Range("L11").Value = 0.509700734558686
$ws.Range("M11").Value = 0.703998080455795
$ws.Range("N11").Value = 0.523760344598141
